$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Step 0: remove the existing _GoBack bookmark (it will be re-added
# later at the end of the new 3rd paragraph).
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------
# Step 1: simple in-place text replacements inside the target
# paragraph (",我所玩的几款" -> ",")
# ---------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("，我所玩的几款", $false, $false, $false, $false, $false, $true, 1, $false, "，", 2)

Write-Host "step1 done"

# ---------------------------------------------------------------
# Step 2: replace the remainder of paragraph 1 with the new text
# for the tail of paragraph 1, plus three brand-new paragraphs
# (paragraph breaks embedded as carriage returns). A placeholder
# token marks where the hyperlink run will be inserted afterwards.
# ---------------------------------------------------------------
$cr = [char]13

$para1Tail = "的游戏发布了新的内容，其中包括我近期玩的暗黑破坏神3和炉石传说，这次发布的内容给我很多启发，在此写下一些感想。"

$para2 = "炉石传说新副本探险者协会中新出了一张宇宙流用回满血的卡雷诺杰克逊，暗黑破坏神2.4补丁中出现了新的散件套戒梦魇者套装。这两者将分别是改变整个游戏环境的两个事物。"

$para3a = "炉石传说自发布以来到现在，期间经历过的各个时期的阶段特点都可以概括为某几套卡组的强势，其中涌现出一系列在之后被削弱的自闭卡组，包括冰法、红龙一刀战、滚石战、龙狗一波猎、奇迹贼、奴隶战。一些非自闭的卡组因为某些卡的强大而被削弱，包括冠军贼、秃鹫放狗猎、送葬亡语猎。同时也有一些卡组因为其优秀的思路以及被暴雪策划所认可的游戏体验而经久不衰，包括咆哮德、动物园，在nga上有篇文章详细说明过过去的几个时期："

$para3b = "。而在如今这个时期，正是奴隶战刚被削弱不久，天梯环境快攻遍地毒瘤横行的黑暗时期。也正是这个时期不久，暴雪出了一张真正能制裁快攻的卡组，与此同时盘活了一种全新的玩法——宇宙流。不管未来环境会不会被人们逐渐厌恶，至少短期来看，这次暴雪做的很成功。"

$para4 = "暗黑破坏神3"

$marker = "@@HYPERLINK@@"

$full = $para1Tail + $cr + $para2 + $cr + $para3a + $marker + $para3b + $cr + $para4

$rng2 = $d.Content
$rng2.Find.Execute("游戏（暗黑破坏神，魔兽世界，炉石传说）也已发布了新的内容，在此写下一些感想。", $false, $false, $false, $false, $false, $true, 1, $false, $full, 2)

Write-Host "step2 done"

# ---------------------------------------------------------------
# Step 3: the new paragraphs 2/3/4 were cloned from paragraph 1's
# paragraph properties (pStyle a3 + numPr + ind firstLineChars=0).
# They must lose the bullet numbering but gain a 420-twip (21pt)
# left indent, matching the target structure. Locate them by their
# unique leading text.
# ---------------------------------------------------------------
function Find-ParaByText($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($needle)) {
            return $p
        }
    }
    return $null
}

$para2 = Find-ParaByText("炉石传说新副本探险者协会中")
$para3 = Find-ParaByText("炉石传说自发布以来")
$para4 = Find-ParaByText("暗黑破坏神3")

foreach ($pp in @($para2, $para3, $para4)) {
    $pp.Range.ListFormat.RemoveNumbers()
    $pp.Format.LeftIndent = 21
    $pp.Format.CharacterUnitFirstLineIndent = 0
}

Write-Host "step3 done"

